$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.535005927085876
$ws.Range("B1").Value = 1.882046580314636
$ws.Range("C1").Value = 2.007616758346558
$ws.Range("D1").Value = 2.318821430206299
$ws.Range("E1").Value = 2.919097900390625
